# Update "想去人数" (want-to-go count) figures, and for two "暂时售罄"
# (temporarily sold out) rows switch the "最低票价" (min price) column
# from a text placeholder to a real numeric price, per the refreshed
# data export.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 (Exhibitions) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1223
$ws1.Range("F4").Value  = 12777
$ws1.Range("F5").Value  = 725
$ws1.Range("F8").Value  = 56
$ws1.Range("F10").Value = 1865
$ws1.Range("F13").Value = 513
$ws1.Range("F15").Value = 124
$ws1.Range("F16").Value = 341
$ws1.Range("F17").Value = 226
$ws1.Range("F18").Value = 294
$ws1.Range("F19").Value = 129
$ws1.Range("F22").Value = 215
$ws1.Range("F23").Value = 243
$ws1.Range("F24").Value = 1280
$ws1.Range("F25").Value = 331
$ws1.Range("F26").Value = 64

# ---- Sheet 2: 演出 (Performances) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value  = 6
$ws2.Range("F6").Value  = 149
$ws2.Range("F7").Value  = 18
$ws2.Range("F8").Value  = 8
$ws2.Range("F9").Value  = 72
$ws2.Range("F10").Value = 72
$ws2.Range("F11").Value = 355
$ws2.Range("F13").Value = 3
$ws2.Range("F20").Value = 11

# ---- Sheet 3: 本地生活 (Local life) ----
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 870
$ws3.Range("F3").Value = 3916
$ws3.Range("G3").Value = 0

# ---- Sheet 4: 全部类型 (All types) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 870
$ws4.Range("F5").Value  = 6
$ws4.Range("F6").Value  = 1223
$ws4.Range("F7").Value  = 12777
$ws4.Range("F9").Value  = 725
$ws4.Range("F10").Value = 3917
$ws4.Range("G10").Value = 30
$ws4.Range("F13").Value = 56
$ws4.Range("F15").Value = 1865
$ws4.Range("F18").Value = 513
$ws4.Range("F21").Value = 149
$ws4.Range("F22").Value = 149
$ws4.Range("F23").Value = 18
$ws4.Range("F24").Value = 124
$ws4.Range("F25").Value = 8
$ws4.Range("F26").Value = 72
$ws4.Range("F27").Value = 72
$ws4.Range("F28").Value = 355
$ws4.Range("F29").Value = 341
$ws4.Range("F31").Value = 226
$ws4.Range("F32").Value = 294
$ws4.Range("F33").Value = 129
$ws4.Range("F36").Value = 3
$ws4.Range("F37").Value = 215
$ws4.Range("F40").Value = 243
$ws4.Range("F41").Value = 1280
$ws4.Range("F43").Value = 331
$ws4.Range("F44").Value = 64
$ws4.Range("F49").Value = 11
